$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 104.833336
$ws.Range("I4").Value = 104.833336
$ws.Range("K4").Value = 104.833336
$ws.Range("M4").Value = 9.166663999999997
$ws.Range("H19").Value = 1282.64
$ws.Range("I19").Value = 982.2
$ws.Range("K19").Value = 982.2
$ws.Range("M19").Value = -807.2
$ws.Range("H33").Value = 409.33334
$ws.Range("I33").Value = 459.4
$ws.Range("K33").Value = 459.4
$ws.Range("M33").Value = -230.4
$ws.Range("H39").Value = 353.23077
$ws.Range("I39").Value = 231.09091
$ws.Range("K39").Value = 693.27273
$ws.Range("M39").Value = -397.27273
$ws.Range("H53").Value = 819.6799999999999
$ws.Range("I53").Value = 653.41174
$ws.Range("K53").Value = 653.41174
$ws.Range("M53").Value = -16.41174000000001
$ws.Range("H61").Value = 4824.6665
$ws.Range("I61").Value = 4824.6665
$ws.Range("K61").Value = 14473.9995
$ws.Range("M61").Value = -14301.9995
$ws.Range("H76").Value = 94506950
$ws.Range("I76").Value = 70382216
$ws.Range("K76").Value = 70382216
$ws.Range("M76").Value = -70381901
$ws.Range("H79").Value = 94506950
$ws.Range("I79").Value = 70382216
$ws.Range("K79").Value = 70382216
$ws.Range("M79").Value = -70381124
$ws.Range("H86").Value = 147799.42
$ws.Range("I86").Value = 203939.4
$ws.Range("K86").Value = 203939.4
$ws.Range("M86").Value = -202816.4
$ws.Range("H89").Value = 147799.42
$ws.Range("I89").Value = 203939.4
$ws.Range("K89").Value = 1019697
$ws.Range("M89").Value = -1014081
$ws.Range("H92").Value = 2589.7632
$ws.Range("I92").Value = 2318.7812
$ws.Range("J92").Value = 4035
$ws.Range("K92").Value = 2318.7812
$ws.Range("L92").Value = 4035
$ws.Range("M92").Value = -1070.7812
$ws.Range("N92").Value = -6531
$ws.Range("H98").Value = 7767152.5
$ws.Range("I98").Value = 8269725.5
$ws.Range("K98").Value = 8269725.5
$ws.Range("M98").Value = -8268227.5
$ws.Range("H100").Value = 2476.8
$ws.Range("I100").Value = 1813
$ws.Range("J100").Value = 3472.5
$ws.Range("K100").Value = 1813
$ws.Range("L100").Value = 3472.5
$ws.Range("M100").Value = -1272
$ws.Range("N100").Value = -4554.5
$ws.Range("H112").Value = 5892480
$ws.Range("J112").Value = 7070796
$ws.Range("L112").Value = 21212388
$ws.Range("N112").Value = -21214604
$ws.Range("H113").Value = 6772.1934
$ws.Range("I113").Value = 7623.2085
$ws.Range("J113").Value = 3854.4285
$ws.Range("K113").Value = 7623.2085
$ws.Range("L113").Value = 3854.4285
$ws.Range("M113").Value = -4369.2085
$ws.Range("N113").Value = -10362.4285
$ws.Range("H122").Value = 7767152.5
$ws.Range("I122").Value = 8269725.5
$ws.Range("K122").Value = 24809176.5
$ws.Range("M122").Value = -24806726.5
$ws.Range("H132").Value = 3926.2407
$ws.Range("I132").Value = 3562.0815
$ws.Range("K132").Value = 10686.2445
$ws.Range("M132").Value = -8156.244499999999
$ws.Range("H137").Value = 1874784.6
$ws.Range("I137").Value = 2594554.8
$ws.Range("J137").Value = 3382.1333
$ws.Range("K137").Value = 7783664.399999999
$ws.Range("L137").Value = 10146.3999
$ws.Range("M137").Value = -7781114.399999999
$ws.Range("N137").Value = -15246.3999
$ws.Range("H138").Value = 2583.75
$ws.Range("I138").Value = 1829.25
$ws.Range("J138").Value = 2877.1667
$ws.Range("K138").Value = 5487.75
$ws.Range("L138").Value = 8631.500100000001
$ws.Range("M138").Value = -347.75
$ws.Range("N138").Value = -18911.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1568.6471
$ws.Range("I2").Value = 1261.2333
$ws.Range("K2").Value = 1261.2333
$ws.Range("M2").Value = -1148.2333
$ws.Range("H4").Value = 205.5
$ws.Range("I4").Value = 105
$ws.Range("K4").Value = 105
$ws.Range("M4").Value = 11
$ws.Range("H32").Value = 1096.69
$ws.Range("I32").Value = 1096.69
$ws.Range("K32").Value = 1096.69
$ws.Range("M32").Value = -809.6900000000001
$ws.Range("H43").Value = 18995.5
$ws.Range("I43").Value = 18995.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 18995.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -18682.5
$ws.Range("N43").ClearContents()
$ws.Range("H45").Value = 1855.4286
$ws.Range("I45").Value = 1797.6
$ws.Range("K45").Value = 1797.6
$ws.Range("M45").Value = -1420.6
$ws.Range("H46").Value = 5674.4
$ws.Range("I46").Value = 4425
$ws.Range("J46").Value = 7548.5
$ws.Range("K46").Value = 4425
$ws.Range("L46").Value = 7548.5
$ws.Range("M46").Value = -4106
$ws.Range("N46").Value = -8186.5
$ws.Range("H61").Value = 2799.9429
$ws.Range("I61").Value = 1854.5
$ws.Range("J61").Value = 4862.727
$ws.Range("K61").Value = 1854.5
$ws.Range("L61").Value = 4862.727
$ws.Range("M61").Value = -1642.5
$ws.Range("N61").Value = -5286.727
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 54953.55
$ws.Range("I74").Value = 63979.938
$ws.Range("J74").Value = 6812.8335
$ws.Range("K74").Value = 63979.938
$ws.Range("L74").Value = 6812.8335
$ws.Range("M74").Value = -63105.938
$ws.Range("N74").Value = -8560.833500000001
$ws.Range("H77").Value = 54953.55
$ws.Range("I77").Value = 63979.938
$ws.Range("J77").Value = 6812.8335
$ws.Range("K77").Value = 319899.69
$ws.Range("L77").Value = 34064.1675
$ws.Range("M77").Value = -315531.69
$ws.Range("N77").Value = -42800.1675
$ws.Range("H110").Value = 2767.1538
$ws.Range("I110").Value = 1402.5555
$ws.Range("J110").Value = 5837.5
$ws.Range("K110").Value = 1402.5555
$ws.Range("L110").Value = 5837.5
$ws.Range("M110").Value = 642.4445000000001
$ws.Range("N110").Value = -9927.5
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H116").Value = 1568.6471
$ws.Range("I116").Value = 1261.2333
$ws.Range("K116").Value = 1261.2333
$ws.Range("M116").Value = 1032.7667
$ws.Range("H122").Value = 2159.5908
$ws.Range("I122").Value = 2159.5908
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6478.7724
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4028.7724
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2609.0557
$ws.Range("I132").Value = 2280.8438
$ws.Range("K132").Value = 6842.5314
$ws.Range("M132").Value = -4312.5314
$ws.Range("H136").Value = 2799.9429
$ws.Range("I136").Value = 1854.5
$ws.Range("J136").Value = 4862.727
$ws.Range("K136").Value = 5563.5
$ws.Range("L136").Value = 14588.181
$ws.Range("M136").Value = -3013.5
$ws.Range("N136").Value = -19688.181
$ws.Range("H139").Value = 69714.5
$ws.Range("J139").Value = 69714.5
$ws.Range("L139").Value = 69714.5
$ws.Range("N139").Value = -79994.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1568.6471
$ws.Range("I3").Value = 1261.2333
$ws.Range("K3").Value = 1261.2333
$ws.Range("M3").Value = -1147.2333
$ws.Range("H20").Value = 8333.5
$ws.Range("I20").Value = 13473.909
$ws.Range("K20").Value = 13473.909
$ws.Range("M20").Value = -13226.909
$ws.Range("H70").Value = 159998
$ws.Range("J70").Value = 159998
$ws.Range("L70").Value = 159998
$ws.Range("N70").Value = -160584
$ws.Range("H73").Value = 159998
$ws.Range("J73").Value = 159998
$ws.Range("L73").Value = 159998
$ws.Range("N73").Value = -162026
$ws.Range("H80").Value = 915.7083
$ws.Range("I80").Value = 1239.1818
$ws.Range("J80").Value = 642
$ws.Range("K80").Value = 1239.1818
$ws.Range("L80").Value = 642
$ws.Range("M80").Value = -241.1818000000001
$ws.Range("N80").Value = -2638
$ws.Range("H83").Value = 915.7083
$ws.Range("I83").Value = 1239.1818
$ws.Range("J83").Value = 642
$ws.Range("K83").Value = 6195.909000000001
$ws.Range("L83").Value = 3210
$ws.Range("M83").Value = -1203.909000000001
$ws.Range("N83").Value = -13194
$ws.Range("H94").Value = 1359.8695
$ws.Range("I94").Value = 1156.2858
$ws.Range("K94").Value = 1156.2858
$ws.Range("M94").Value = -705.2858000000001
$ws.Range("H105").Value = 2022.625
$ws.Range("I105").Value = 1980.1305
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 1980.1305
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -233.1305
$ws.Range("N105").Value = -6494
$ws.Range("H134").Value = 4800.15
$ws.Range("I134").Value = 3751
$ws.Range("J134").Value = 6373.875
$ws.Range("K134").Value = 11253
$ws.Range("L134").Value = 19121.625
$ws.Range("M134").Value = -8718
$ws.Range("N134").Value = -24191.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H16").Value = 1801.6818
$ws.Range("I16").Value = 1770.8
$ws.Range("K16").Value = 1770.8
$ws.Range("M16").Value = -1483.8
$ws.Range("H31").Value = 183358.86
$ws.Range("I31").Value = 228408.1
$ws.Range("J31").Value = 3161.9092
$ws.Range("K31").Value = 228408.1
$ws.Range("L31").Value = 3161.9092
$ws.Range("M31").Value = -228113.1
$ws.Range("N31").Value = -3751.9092
$ws.Range("H34").Value = 183358.86
$ws.Range("I34").Value = 228408.1
$ws.Range("J34").Value = 3161.9092
$ws.Range("K34").Value = 228408.1
$ws.Range("L34").Value = 3161.9092
$ws.Range("M34").Value = -228206.1
$ws.Range("N34").Value = -3565.9092
$ws.Range("H58").Value = 2264.2856
$ws.Range("I58").Value = 2264.2856
$ws.Range("K58").Value = 2264.2856
$ws.Range("M58").Value = -2061.2856
$ws.Range("H105").Value = 5979
$ws.Range("I105").Value = 1982.375
$ws.Range("K105").Value = 1982.375
$ws.Range("M105").Value = -235.375
$ws.Range("H107").Value = 3937.9473
$ws.Range("I107").Value = 228.375
$ws.Range("J107").Value = 6635.8184
$ws.Range("K107").Value = 228.375
$ws.Range("L107").Value = 6635.8184
$ws.Range("M107").Value = 1691.625
$ws.Range("N107").Value = -10475.8184
$ws.Range("H113").Value = 1801.6818
$ws.Range("I113").Value = 1770.8
$ws.Range("K113").Value = 1770.8
$ws.Range("M113").Value = 399.2
$ws.Range("H122").Value = 4947.5
$ws.Range("I122").Value = 4947.5
$ws.Range("K122").Value = 14842.5
$ws.Range("M122").Value = -12392.5
$ws.Range("H132").Value = 1775.5135
$ws.Range("I132").Value = 1822.8572
$ws.Range("K132").Value = 5468.571599999999
$ws.Range("M132").Value = -2938.571599999999
$ws.Range("H134").Value = 5528.1816
$ws.Range("I134").Value = 6285.1055
$ws.Range("K134").Value = 18855.3165
$ws.Range("M134").Value = -16320.3165
$ws.Range("H136").Value = 2264.2856
$ws.Range("I136").Value = 2264.2856
$ws.Range("K136").Value = 6792.8568
$ws.Range("M136").Value = -4242.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1345.1177
$ws.Range("I5").Value = 623.9091
$ws.Range("J5").Value = 2667.3333
$ws.Range("K5").Value = 1871.7273
$ws.Range("L5").Value = 8001.999899999999
$ws.Range("M5").Value = -1759.7273
$ws.Range("N5").Value = -8225.999899999999
$ws.Range("H29").Value = 4599.1665
$ws.Range("J29").Value = 4773
$ws.Range("L29").Value = 14319
$ws.Range("N29").Value = -14873
$ws.Range("H39").Value = 1255.5416
$ws.Range("I39").Value = 883.3333
$ws.Range("J39").Value = 1478.8667
$ws.Range("K39").Value = 2649.9999
$ws.Range("L39").Value = 4436.6001
$ws.Range("M39").Value = -2355.9999
$ws.Range("N39").Value = -5024.6001
$ws.Range("H40").Value = 89.833336
$ws.Range("I40").Value = 84.75
$ws.Range("K40").Value = 339
$ws.Range("M40").Value = -270
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 6000
$ws.Range("M62").Value = -5314
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 18000
$ws.Range("M65").Value = -14568
$ws.Range("H69").Value = 750
$ws.Range("I69").Value = 750
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 2250
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -1439
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 750
$ws.Range("I72").Value = 750
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 6750
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -2694
$ws.Range("N72").ClearContents()
$ws.Range("H92").Value = 292.85715
$ws.Range("I92").Value = 239.5
$ws.Range("K92").Value = 718.5
$ws.Range("M92").Value = 529.5
$ws.Range("H113").Value = 855.41174
$ws.Range("J113").Value = 907.1667
$ws.Range("L113").Value = 2721.5001
$ws.Range("N113").Value = -7061.5001
$ws.Range("H114").Value = 998.5
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("H117").Value = 890.25
$ws.Range("J117").Value = 1250
$ws.Range("L117").Value = 3750
$ws.Range("N117").Value = -10634
$ws.Range("H121").Value = 79698.28999999999
$ws.Range("I121").Value = 6996.5
$ws.Range("J121").Value = 91815.25
$ws.Range("K121").Value = 20989.5
$ws.Range("L121").Value = 275445.75
$ws.Range("M121").Value = -19679.5
$ws.Range("N121").Value = -278065.75
$ws.Range("H131").Value = 8475953
$ws.Range("I131").Value = 62500868
$ws.Range("K131").Value = 187502604
$ws.Range("M131").Value = -187497564
$ws.Range("H132").Value = 7154.35
$ws.Range("I132").Value = 7693.222
$ws.Range("J132").Value = 2304.5
$ws.Range("K132").Value = 69238.99799999999
$ws.Range("L132").Value = 20740.5
$ws.Range("M132").Value = -66708.99799999999
$ws.Range("N132").Value = -25800.5
$ws.Range("H135").Value = 1345.1177
$ws.Range("I135").Value = 623.9091
$ws.Range("J135").Value = 2667.3333
$ws.Range("K135").Value = 5615.1819
$ws.Range("L135").Value = 24005.9997
$ws.Range("M135").Value = -3080.1819
$ws.Range("N135").Value = -29075.9997
$ws.Range("H140").Value = 5683548.5
$ws.Range("I140").Value = 13889873
$ws.Range("J140").Value = 2246.5386
$ws.Range("K140").Value = 41669619
$ws.Range("L140").Value = 6739.6158
$ws.Range("M140").Value = -41664439
$ws.Range("N140").Value = -17099.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 141.5
$ws.Range("I2").Value = 159.33333
$ws.Range("J2").Value = 88
$ws.Range("K2").Value = 159.33333
$ws.Range("L2").Value = 88
$ws.Range("M2").Value = -46.33332999999999
$ws.Range("N2").Value = -314
$ws.Range("H80").Value = 4398.706
$ws.Range("I80").Value = 3222
$ws.Range("J80").Value = 4760.769
$ws.Range("K80").Value = 3222
$ws.Range("L80").Value = 4760.769
$ws.Range("M80").Value = -2224
$ws.Range("N80").Value = -6756.769
$ws.Range("H83").Value = 4398.706
$ws.Range("I83").Value = 3222
$ws.Range("J83").Value = 4760.769
$ws.Range("K83").Value = 16110
$ws.Range("L83").Value = 23803.845
$ws.Range("M83").Value = -11118
$ws.Range("N83").Value = -33787.845
$ws.Range("H97").Value = 1982.6
$ws.Range("I97").Value = 1540.5667
$ws.Range("J97").Value = 2866.6667
$ws.Range("K97").Value = 1540.5667
$ws.Range("L97").Value = 2866.6667
$ws.Range("M97").Value = -1044.5667
$ws.Range("N97").Value = -3858.6667
$ws.Range("H102").Value = 2607.6553
$ws.Range("I102").Value = 2546.8462
$ws.Range("K102").Value = 2546.8462
$ws.Range("M102").Value = -924.8462
$ws.Range("H107").Value = 479.33334
$ws.Range("I107").Value = 550.4
$ws.Range("K107").Value = 550.4
$ws.Range("M107").Value = 1369.6
$ws.Range("H111").Value = 77665
$ws.Range("J111").Value = 77665
$ws.Range("L111").Value = 77665
$ws.Range("N111").Value = -83799
$ws.Range("H126").Value = 6178.3335
$ws.Range("I126").Value = 8624.429
$ws.Range("J126").Value = 4038
$ws.Range("K126").Value = 25873.287
$ws.Range("L126").Value = 12114
$ws.Range("M126").Value = -23403.287
$ws.Range("N126").Value = -17054
$ws.Range("H132").Value = 20551.15
$ws.Range("I132").Value = 24319
$ws.Range("J132").Value = 7128.1875
$ws.Range("K132").Value = 72957
$ws.Range("L132").Value = 21384.5625
$ws.Range("M132").Value = -70427
$ws.Range("N132").Value = -26444.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8103.44
$ws.Range("I7").Value = 8285.956
$ws.Range("J7").Value = 6004.5
$ws.Range("K7").Value = 8285.956
$ws.Range("L7").Value = 6004.5
$ws.Range("M7").Value = -8173.956
$ws.Range("N7").Value = -6228.5
$ws.Range("H17").Value = 430
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H22").Value = 4062.6316
$ws.Range("I22").Value = 4135.0713
$ws.Range("K22").Value = 4135.0713
$ws.Range("M22").Value = -3840.0713
$ws.Range("H27").Value = 4062.6316
$ws.Range("I27").Value = 4135.0713
$ws.Range("K27").Value = 4135.0713
$ws.Range("M27").Value = -4028.0713
$ws.Range("H40").Value = 21109.416
$ws.Range("I40").Value = 29416.5
$ws.Range("J40").Value = 4495.25
$ws.Range("K40").Value = 29416.5
$ws.Range("L40").Value = 4495.25
$ws.Range("M40").Value = -29280.5
$ws.Range("N40").Value = -4767.25
$ws.Range("H46").Value = 3485.0715
$ws.Range("I46").Value = 1643.6666
$ws.Range("J46").Value = 6799.6
$ws.Range("K46").Value = 1643.6666
$ws.Range("L46").Value = 6799.6
$ws.Range("M46").Value = -1455.6666
$ws.Range("N46").Value = -7175.6
$ws.Range("H55").Value = 7299.4346
$ws.Range("I55").Value = 760.44446
$ws.Range("K55").Value = 760.44446
$ws.Range("M55").Value = -587.44446
$ws.Range("H61").Value = 1592.875
$ws.Range("I61").Value = 1787
$ws.Range("J61").Value = 1398.75
$ws.Range("K61").Value = 1787
$ws.Range("L61").Value = 1398.75
$ws.Range("M61").Value = -1585
$ws.Range("N61").Value = -1802.75
$ws.Range("H68").Value = 3817.4
$ws.Range("I68").Value = 3817.4
$ws.Range("K68").Value = 3817.4
$ws.Range("M68").Value = -3068.4
$ws.Range("H71").Value = 3817.4
$ws.Range("I71").Value = 3817.4
$ws.Range("K71").Value = 19087
$ws.Range("M71").Value = -15343
$ws.Range("H93").Value = 66667760
$ws.Range("I93").Value = 664.375
$ws.Range("J93").Value = 142858740
$ws.Range("K93").Value = 664.375
$ws.Range("L93").Value = 142858740
$ws.Range("M93").Value = 583.625
$ws.Range("N93").Value = -142861236
$ws.Range("H100").Value = 2607.1667
$ws.Range("I100").Value = 2607.1667
$ws.Range("K100").Value = 2607.1667
$ws.Range("M100").Value = -2066.1667
$ws.Range("H113").Value = 1592.875
$ws.Range("I113").Value = 1787
$ws.Range("J113").Value = 1398.75
$ws.Range("K113").Value = 1787
$ws.Range("L113").Value = 1398.75
$ws.Range("M113").Value = 383
$ws.Range("N113").Value = -5738.75
$ws.Range("H122").Value = 5820.154
$ws.Range("I122").Value = 5646.727
$ws.Range("J122").Value = 6774
$ws.Range("K122").Value = 16940.181
$ws.Range("L122").Value = 20322
$ws.Range("M122").Value = -14490.181
$ws.Range("N122").Value = -25222
$ws.Range("H126").Value = 8103.44
$ws.Range("I126").Value = 8285.956
$ws.Range("J126").Value = 6004.5
$ws.Range("K126").Value = 24857.868
$ws.Range("L126").Value = 18013.5
$ws.Range("M126").Value = -22387.868
$ws.Range("N126").Value = -22953.5
$ws.Range("H132").Value = 4023.554
$ws.Range("I132").Value = 2824.361
$ws.Range("J132").Value = 5512.207
$ws.Range("K132").Value = 8473.082999999999
$ws.Range("L132").Value = 16536.621
$ws.Range("M132").Value = -5943.082999999999
$ws.Range("N132").Value = -21596.621
$ws.Range("H136").Value = 1962.8379
$ws.Range("I136").Value = 1523.8064
$ws.Range("J136").Value = 4231.1665
$ws.Range("K136").Value = 4571.4192
$ws.Range("L136").Value = 12693.4995
$ws.Range("M136").Value = -2021.4192
$ws.Range("N136").Value = -17793.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6257.231
$ws.Range("J14").Value = 6528.75
$ws.Range("L14").Value = 6528.75
$ws.Range("N14").Value = -6864.75
$ws.Range("H82").Value = 68000
$ws.Range("J82").Value = 68000
$ws.Range("L82").Value = 68000
$ws.Range("N82").Value = -68766
$ws.Range("H85").Value = 68000
$ws.Range("J85").Value = 68000
$ws.Range("L85").Value = 68000
$ws.Range("N85").Value = -70652
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H100").Value = 1035.875
$ws.Range("I100").Value = 1035.875
$ws.Range("K100").Value = 2071.75
$ws.Range("M100").Value = -1530.75
$ws.Range("H107").Value = 768.54285
$ws.Range("I107").Value = 707.2414
$ws.Range("K107").Value = 2121.7242
$ws.Range("M107").Value = -201.7242000000001
$ws.Range("H126").Value = 1917.4
$ws.Range("I126").Value = 1977.5
$ws.Range("J126").Value = 1777.1666
$ws.Range("K126").Value = 5932.5
$ws.Range("L126").Value = 5331.4998
$ws.Range("M126").Value = -3462.5
$ws.Range("N126").Value = -10271.4998
$ws.Range("H132").Value = 1320.1143
$ws.Range("I132").Value = 779.3333
$ws.Range("K132").Value = 2337.9999
$ws.Range("M132").Value = 192.0001000000002
$ws.Range("H136").Value = 170749.02
$ws.Range("I136").Value = 200293.05
$ws.Range("J136").Value = 3332.889
$ws.Range("K136").Value = 600879.1499999999
$ws.Range("L136").Value = 9998.667000000001
$ws.Range("M136").Value = -598329.1499999999
$ws.Range("N136").Value = -15098.667
$ws.Range("H139").Value = 66666.664
$ws.Range("J139").Value = 66666.664
$ws.Range("L139").Value = 66666.664
$ws.Range("N139").Value = -76946.664
